$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value = 16.18989999999999
$ws.Range("C7").Value = -12.6809
$ws.Range("A9").Value = -21.73770000000002
$ws.Range("C12").Value = -10.9187
$ws.Range("C14").Value = -12.4282
$ws.Range("E15").Value = 16.471
$ws.Range("A18").Value = -22.14340000000001
$ws.Range("A20").Value = -21.68039999999998
$ws.Range("C26").Value = -12.5658
$ws.Range("A27").Value = -21.8764
$ws.Range("C27").Value = -12.7798
$ws.Range("C29").Value = -11.10680000000001
$ws.Range("E33").Value = 17.29260000000001
$ws.Range("A35").Value = -21.60459999999998
$ws.Range("E35").Value = 16.7474
$ws.Range("C37").Value = -13.1597
$ws.Range("C38").Value = -12.6449
$ws.Range("E38").Value = 16.60289999999999
$ws.Range("E43").Value = 17.26930000000002
$ws.Range("E44").Value = 16.6899
$ws.Range("E47").Value = 16.43249999999999
$ws.Range("C51").Value = -11.6022
$ws.Range("E51").Value = 17.21620000000001
$ws.Range("C52").Value = -11.3044
$ws.Range("C55").Value = -13.58559999999999
$ws.Range("E57").Value = 16.44090000000001
$ws.Range("E63").Value = 18.47820000000001
$ws.Range("A69").Value = -21.711
$ws.Range("C69").Value = -11.0936
$ws.Range("C70").Value = -11.8943
$ws.Range("E70").Value = 17.3769
$ws.Range("A76").Value = -20.1432
$ws.Range("A78").Value = -20.34579999999999
$ws.Range("C81").Value = -12.5503
$ws.Range("A82").Value = -21.9006
$ws.Range("A83").Value = -22.0501
$ws.Range("C83").Value = -12.91919999999999
$ws.Range("E88").Value = 16.4405
$ws.Range("A93").Value = -21.06579999999999
$ws.Range("E99").Value = 16.5032
$ws.Range("C102").Value = -13.3717
